$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every touched cell to remain plain text so values such as
# "305.52", "-4.05%", "21" round-trip byte-for-byte instead of being
# coerced into numbers/percentages/dates by COM type inference.
$cells = @(
    "D2", "E2", "G2", "E3", "G3", "D4", "E4", "G4", "D5", "G5", "D6", "E6",
    "G6", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10",
    "D11", "E11", "G11", "D12", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14",
    "D15", "E15", "G15", "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18",
    "D19", "E19", "G19", "D20", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22",
    "D23", "E23", "G23", "D24", "E24", "G24", "D25", "E25", "G25", "D26", "E26", "G26",
    "D27", "E27", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36",
    "G37", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42",
    "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46",
    "E46", "G46", "D47", "E47", "G47", "B48", "C48", "D48", "E48", "G48", "B49", "C49",
    "D49", "E49", "G49", "D50", "E50", "G50", "D51", "E51", "G51"
)
foreach ($ref in $cells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values (as text).
$ws.Range("D2").Value = '305.52'
$ws.Range("E2").Value = '-4.05%'
$ws.Range("G2").Value = '21'
$ws.Range("E3").Value = '-6.54%'
$ws.Range("G3").Value = '21'
$ws.Range("D4").Value = '5.097'
$ws.Range("E4").Value = '-0.91%'
$ws.Range("G4").Value = '21'
$ws.Range("D5").Value = '0.07699'
$ws.Range("G5").Value = '21'
$ws.Range("D6").Value = '4.359'
$ws.Range("E6").Value = '0.20%'
$ws.Range("G6").Value = '21'
$ws.Range("E7").Value = '-1.69%'
$ws.Range("G7").Value = '21'
$ws.Range("D8").Value = '1.868'
$ws.Range("E8").Value = '-8.54%'
$ws.Range("G8").Value = '21'
$ws.Range("D9").Value = '3.194'
$ws.Range("E9").Value = '-2.96%'
$ws.Range("G9").Value = '21'
$ws.Range("D10").Value = '0.9174'
$ws.Range("E10").Value = '-2.37%'
$ws.Range("G10").Value = '21'
$ws.Range("D11").Value = '0.1172'
$ws.Range("E11").Value = '-13.50%'
$ws.Range("G11").Value = '21'
$ws.Range("D12").Value = '0.1876'
$ws.Range("E12").Value = '-4.93%'
$ws.Range("G12").Value = '21'
$ws.Range("D13").Value = '0.08732'
$ws.Range("E13").Value = '-4.35%'
$ws.Range("G13").Value = '21'
$ws.Range("D14").Value = '0.03427'
$ws.Range("E14").Value = '-2.42%'
$ws.Range("G14").Value = '21'
$ws.Range("D15").Value = '0.09692'
$ws.Range("E15").Value = '-1.10%'
$ws.Range("G15").Value = '21'
$ws.Range("D16").Value = '0.001366'
$ws.Range("E16").Value = '-3.18%'
$ws.Range("G16").Value = '21'
$ws.Range("D17").Value = '0.005922'
$ws.Range("E17").Value = '-2.79%'
$ws.Range("G17").Value = '21'
$ws.Range("D18").Value = '3.566'
$ws.Range("E18").Value = '-3.38%'
$ws.Range("G18").Value = '21'
$ws.Range("D19").Value = '0.3373'
$ws.Range("E19").Value = '-3.11%'
$ws.Range("G19").Value = '21'
$ws.Range("D20").Value = '0.1277'
$ws.Range("E20").Value = '-3.06%'
$ws.Range("G20").Value = '21'
$ws.Range("D21").Value = '5.019'
$ws.Range("E21").Value = '1.44%'
$ws.Range("G21").Value = '21'
$ws.Range("D22").Value = '0.2499'
$ws.Range("E22").Value = '1.66%'
$ws.Range("G22").Value = '21'
$ws.Range("D23").Value = '0.02112'
$ws.Range("E23").Value = '5,164.77%'
$ws.Range("G23").Value = '21'
$ws.Range("D24").Value = '0.04325'
$ws.Range("E24").Value = '-1.05%'
$ws.Range("G24").Value = '21'
$ws.Range("D25").Value = '0.001215'
$ws.Range("E25").Value = '-1.32%'
$ws.Range("G25").Value = '21'
$ws.Range("D26").Value = '0.004473'
$ws.Range("E26").Value = '-6.75%'
$ws.Range("G26").Value = '21'
$ws.Range("D27").Value = '0.0001353'
$ws.Range("E27").Value = '3.95%'
$ws.Range("G27").Value = '21'
$ws.Range("G28").Value = '21'
$ws.Range("G29").Value = '21'
$ws.Range("G30").Value = '21'
$ws.Range("G31").Value = '21'
$ws.Range("G32").Value = '21'
$ws.Range("G33").Value = '21'
$ws.Range("G34").Value = '21'
$ws.Range("G35").Value = '21'
$ws.Range("G36").Value = '21'
$ws.Range("G37").Value = '21'
$ws.Range("G38").Value = '21'
$ws.Range("D39").Value = '0.02226'
$ws.Range("E39").Value = '-0.81%'
$ws.Range("G39").Value = '21'
$ws.Range("D40").Value = '0.04899'
$ws.Range("E40").Value = '-5.93%'
$ws.Range("G40").Value = '21'
$ws.Range("D41").Value = '0.007558'
$ws.Range("E41").Value = '-2.55%'
$ws.Range("G41").Value = '21'
$ws.Range("D42").Value = '0.009846'
$ws.Range("E42").Value = '1.21%'
$ws.Range("G42").Value = '21'
$ws.Range("D43").Value = '0.1331'
$ws.Range("E43").Value = '-5.04%'
$ws.Range("G43").Value = '21'
$ws.Range("D44").Value = '0.001998'
$ws.Range("E44").Value = '-2.50%'
$ws.Range("G44").Value = '21'
$ws.Range("D45").Value = '0.008811'
$ws.Range("E45").Value = '-4.17%'
$ws.Range("G45").Value = '21'
$ws.Range("D46").Value = '0.00006687'
$ws.Range("E46").Value = '1.93%'
$ws.Range("G46").Value = '21'
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").Value = '-0.03%'
$ws.Range("G47").Value = '21'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = '0.003005'
$ws.Range("E48").Value = '1.86%'
$ws.Range("G48").Value = '21'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").Value = '0.001303'
$ws.Range("E49").Value = '-23.09%'
$ws.Range("G49").Value = '21'
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").Value = '-0.03%'
$ws.Range("G50").Value = '21'
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").Value = '-0.03%'
$ws.Range("G51").Value = '21'
